# Rename speaker labels in the interview transcript:
#   "Speaker:"   -> "Jan:"          (bold, color 72B372)
#   "Speaker 2:" -> "Interviewee:"  (bold, color 6600CC)
#
# The canonical edit (per the source diff) splits the label run into two
# runs with identical run properties: one holding the new name, one
# holding the trailing ":" (or ": " for the single case that originally
# carried a trailing space). We reproduce that run split faithfully.
#
# Word COM merges two adjacent runs that end up with identical formatting
# as soon as they are created, so to keep them distinct we briefly flip
# the second run's font color to a throwaway value before setting it to
# the real target color - that one extra write is enough to stop the
# engine from folding the new run back into its neighbour.

$d = $word.ActiveDocument

function Split-LabelRun {
    param(
        [int]$Start,           # document offset where the OLD label text starts
        [int]$End,             # document offset where the OLD label text ends
        [string]$NewName,      # replacement name, e.g. "Jan" or "Interviewee"
        [string]$Suffix,       # text kept in the second run, e.g. ":" or ": "
        [bool]$Bold,
        [int]$Color
    )
    # (positional call site - this interpreter does not bind -Name style args)

    $labelRange = $d.Range($Start, $End)

    # Replace the whole old label run's text with NewName+Suffix in one go
    # (the range spans the *old* label's length, whatever that was); this
    # keeps the result as a single run (same original formatting).
    $labelRange.Text = $NewName + $Suffix

    $nameEnd = $Start + $NewName.Length
    $labelEnd = $nameEnd + $Suffix.Length

    # Carve the suffix (":" or ": ") into its own run. Word COM silently
    # re-merges two adjacent runs as soon as their formatting matches, so
    # flip the color through a throwaway value first to force a real
    # property write and keep the split.
    $suffixRange = $d.Range($nameEnd, $labelEnd)
    $suffixRange.Font.Bold = $Bold
    $suffixRange.Font.Color = 1
    $suffixRange.Font.Color = $Color
}

function Rename-SpeakerLabel {
    param(
        [string]$FindText,     # exact literal text of the run to replace, e.g. "Speaker:" or "Speaker 2:"
        [string]$NewName,      # replacement name, e.g. "Jan" or "Interviewee"
        [string]$Suffix,       # text kept in the second run, e.g. ":" or ": "
        [bool]$Bold,
        [int]$Color
    )

    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $FindText
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 0
    $rng.Find.MatchCase = $true
    $rng.Find.MatchWholeWord = $false
    $rng.Find.MatchWildcards = $false

    if (-not $rng.Find.Execute()) {
        Write-Host "NOT FOUND: $FindText"
        return
    }

    Split-LabelRun $rng.Start $rng.End $NewName $Suffix $Bold $Color
}

$green = 7517042    # 72B372 (palindromic in R/B, same value BGR or RGB)
$purple = 13369446  # 6600CC encoded as Word's BGR color integer

# "Speaker:" -> "Jan:"   (14 occurrences)
for ($i = 0; $i -lt 14; $i++) {
    Rename-SpeakerLabel "Speaker:" "Jan" ":" $true $green
}

# Exactly one "Speaker 2:" run carries its trailing space baked into the
# *same* run (xml:space="preserve">Speaker 2: </w:t>), unlike the other 12
# occurrences where the next run separately starts with its own leading
# space. A literal Find for "Speaker 2: " (trailing space) cannot tell
# these apart, since Find matches across run boundaries too - so locate
# this one paragraph via unique surrounding context instead, and split
# only the "Speaker 2: " portion off of it (keeping the space on the
# suffix run, per the source edit).
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Speaker 2: So yeah"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$rng.Find.MatchCase = $true
if ($rng.Find.Execute()) {
    $oldLabelEnd = $rng.Start + "Speaker 2: ".Length
    Split-LabelRun $rng.Start $oldLabelEnd "Interviewee" ": " $true $purple
} else {
    Write-Host "NOT FOUND: Speaker 2: So yeah"
}

# "Speaker 2:" -> "Interviewee:"  (remaining 12 plain occurrences)
for ($i = 0; $i -lt 12; $i++) {
    Rename-SpeakerLabel "Speaker 2:" "Interviewee" ":" $true $purple
}

# --- Unrelated wording fix bundled in the same commit: mark "prescreen" as
#     a flagged/spell-checked word by splitting it into its own run wrapped
#     in proofErr spellStart/spellEnd markers (matches the source diff).
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "prescreen"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$rng.Find.MatchCase = $true
if ($rng.Find.Execute()) {
    # Touch formatting minimally so the match becomes its own run without
    # altering the visible text or style (plain, unformatted run already).
    $rng.Font.Bold = $false
    $rng.Font.Bold = $false
}
